$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.476.47"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "3.195.21"
$ws.Range("E3").Value = "  -3.37%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.39"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.56"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "3.191.42"
$ws.Range("E8").Value = "  -3.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  -4.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.28"
$ws.Range("E11").Value = "  -4.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.69"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "3.716.53"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "3.188.37"
$ws.Range("E17").Value = "  -3.62%  "
$ws.Range("D18").Value = "63.341.68"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.60"
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.92"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.03"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.702"
$ws.Range("E22").Value = "  -4.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.67"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.46"
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.86"
$ws.Range("E25").Value = "  -2.82%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.81"
$ws.Range("E29").Value = "  -3.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.81"
$ws.Range("E30").Value = "  -4.77%  "
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.47"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.40"
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.04"
$ws.Range("E35").Value = "  -4.62%  "
$ws.Range("E36").Value = "  -2.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.42"
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("D38").Value = "0.0₃0726"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0390"
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "405.58"
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.13"
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("E43").Value = "  -6.00%  "
$ws.Range("D44").Value = "2.819.84"
$ws.Range("E44").Value = "  -9.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.255"
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.10"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.70"
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.49"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("E51").Value = "  -1.23%  "

Write-Host "Updated cryptos list"
